$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: new entry for 2024-02-14
$ws.Range("A13").Value = 20240214
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 6
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 5

# Row 14: only G14 has a value
$ws.Range("G14").Value = 8

# Row 15: new entry for 2024-02-29
$ws.Range("A15").Value = 20240229
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 6

# Move the active selection to A13, matching the saved workbook state
$ws.Range("A13").Select() | Out-Null
